$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "privacy values clean" right after "Services
#    clean" (so the sheet order becomes: Services clean, privacy values
#    clean, Services test, Privacy values).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("Services clean")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "privacy values clean"

# Worksheet references taken before the Add() call can go stale (positional
# rebind), so re-acquire both sheets by name now that the sheet collection is
# stable.
$src = $wb.Worksheets.Item("Privacy values")
$dst = $wb.Worksheets.Item("privacy values clean")

# Copy the small summary table (A1:S8) from "Privacy values" into the new
# sheet - values first, then formats (keeps header cell styles s="1"/s="2").
$src.Range("A1:S8").Copy()
$dst.Range("A1").PasteSpecial(-4104)
$src.Range("A1:S8").Copy()
$dst.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Swap the ranges the two hidden chart-backing defined names point to.
# ---------------------------------------------------------------------------
$name0 = $wb.Names.Item("_xlchart.v1.0")
$name1 = $wb.Names.Item("_xlchart.v1.1")
$name0.RefersTo = "='Services test'!`$AE`$2:`$AE`$76"
$name1.RefersTo = "='Services test'!`$AD`$2:`$AD`$76"

# ---------------------------------------------------------------------------
# 3. Restore per-sheet selections / scroll positions.
# ---------------------------------------------------------------------------
$servicesClean = $wb.Worksheets.Item("Services clean")
$servicesClean.Activate()
$servicesClean.Range("AB2").Select()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 6

$privacyValues = $wb.Worksheets.Item("Privacy values")
$privacyValues.Activate()
$privacyValues.Range("A1:S8").Select()

# ---------------------------------------------------------------------------
# 4. Leave "privacy values clean" as the active tab/selection, matching the
#    recorded activeTab of the saved workbook.
# ---------------------------------------------------------------------------
$dst = $wb.Worksheets.Item("privacy values clean")
$dst.Activate()
$dst.Range("B2").Select()

# ---------------------------------------------------------------------------
# 5. Best-effort restore of the workbook window position.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Left = 28680
$excel.ActiveWindow.Top = -120
